$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Cells in columns B/C/D are forced to Text format before assignment so that
# numeric-looking strings (e.g. "583.71") are preserved exactly as text rather
# than being auto-converted into floating point numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.269.00'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.601.57'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.71'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '190.01'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.595.15'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.665'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.08'
$ws.Range("E12").Value = '  -4.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000313'
$ws.Range("E13").Value = '  +7.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.73'
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.181.47'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.601.21'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.265.93'
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '490.49'
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '20.20'
$ws.Range("E23").Value = '  +3.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.96'
$ws.Range("E24").Value = '  -7.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.59'
$ws.Range("E25").Value = '  +7.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.40'
$ws.Range("E26").Value = '  -1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.99'
$ws.Range("E27").Value = '  -4.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.14'
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.54'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.44'
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.64'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.28'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.119'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.39'
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '579.91'
$ws.Range("E35").Value = '  -7.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.08'
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.399'
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("E40").Value = '  +20.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("E41").Value = '  +7.01%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  -6.00%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.45'
$ws.Range("E43").Value = '  -2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.222.54'
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.08'
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0448'
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("E47").Value = '  +5.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.35'
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("E51").Value = '  -2.24%  '
